$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("deals")

# ---- Header row (row 1) ----
$headers = @(
    "title",
    "company",
    "primarycontact",
    "amount",
    "probability",
    "commission",
    "identifier",
    "tags",
    "descrption",
    "nextStep",
    "product",
    "quantity",
    "type",
    "source"
)
for ($c = 1; $c -le $headers.Length; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $cell.Value = $headers[$c - 1]
    $cell.Interior.ColorIndex = 6
}

# ---- Data row 2 ----
$ws.Cells.Item(2, 1).Value = "deals title - 1"
$ws.Cells.Item(2, 2).Value = "Flipkart company"
$ws.Cells.Item(2, 3).Value = "bansal - contact"
$ws.Cells.Item(2, 4).Value = 5000
$ws.Cells.Item(2, 5).Value = 80
$ws.Cells.Item(2, 6).Value = 20
$ws.Cells.Item(2, 7).Value = "test identifier"
$ws.Cells.Item(2, 8).Value = "tagOne, tagTwo, tagThree"
$ws.Cells.Item(2, 9).Value = "test desc -- added by salesperson"
$ws.Cells.Item(2, 10).Value = "waiting for answer from client"
$ws.Cells.Item(2, 11).Value = "Test Product"
$ws.Cells.Item(2, 12).Value = 3
$ws.Cells.Item(2, 13).Value = "Priority"
$ws.Cells.Item(2, 14).Value = "Word of Mouth"

# ---- Data row 3 ----
$ws.Cells.Item(3, 1).Value = "deals title - 2"
$ws.Cells.Item(3, 2).Value = "Amazon Ccompany"
$ws.Cells.Item(3, 3).Value = "John - Contact"
$ws.Cells.Item(3, 4).Value = 8000
$ws.Cells.Item(3, 5).Value = 60
$ws.Cells.Item(3, 6).Value = 10
$ws.Cells.Item(3, 7).Value = "Test - 2"
$ws.Cells.Item(3, 8).Value = "tagFour, tagFive"
$ws.Cells.Item(3, 9).Value = "test desc -- added by salesperson (amazon)"
$ws.Cells.Item(3, 10).Value = "amazon - next step"
$ws.Cells.Item(3, 11).Value = "Test Product"
$ws.Cells.Item(3, 12).Value = 5
$ws.Cells.Item(3, 13).Value = "Priority"
$ws.Cells.Item(3, 14).Value = "Word of Mouth"

# ---- Column widths to roughly match the authored layout ----
$widths = @{
    1 = 12.43
    2 = 18
    3 = 14.86
    5 = 10.71
    6 = 11.57
    7 = 13.43
    8 = 24.14
    9 = 40.71
    10 = 28.14
    11 = 12
    12 = 8.43
    13 = 7.57
    14 = 14.57
}
foreach ($col in $widths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $widths[$col]
}

# ---- Activate deals sheet, make it the selected tab, and set the selection ----
$ws.Activate()
$ws.Range("M11").Select()
